# Update the team-specific transition-probability matrix values on Sheet1
# (rows 2-13 and 15-19) to match the refreshed "team spec time" computation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1766467065868264
$ws.Range("C2").Value = 0.5838323353293413
$ws.Range("J2").Value = 0.008982035928143712
$ws.Range("P2").Value = 0.1437125748502994
$ws.Range("S2").Value = 0.08682634730538923
$ws.Range("B3").Value = 0.004950495049504951
$ws.Range("C3").Value = 0.0297029702970297
$ws.Range("J3").Value = 0.009900990099009901
$ws.Range("P3").Value = 0.7673267326732673
$ws.Range("S3").Value = 0.1881188118811881
$ws.Range("J4").Value = 0.06
$ws.Range("P4").Value = 0.48
$ws.Range("S4").Value = 0.46
$ws.Range("P5").Value = 1
$ws.Range("B6").Value = 0.05761316872427984
$ws.Range("D6").Value = 0.01646090534979424
$ws.Range("E6").Value = 0.00411522633744856
$ws.Range("F6").Value = 0.04526748971193416
$ws.Range("J6").Value = 0.3127572016460906
$ws.Range("O6").Value = 0.04526748971193416
$ws.Range("Q6").Value = 0.1193415637860082
$ws.Range("R6").Value = 0.07818930041152264
$ws.Range("S6").Value = 0.3209876543209876
$ws.Range("B7").Value = 0.1353711790393013
$ws.Range("D7").Value = 0.03493449781659388
$ws.Range("F7").Value = 0.05676855895196507
$ws.Range("J7").Value = 0.09606986899563319
$ws.Range("O7").Value = 0.02620087336244541
$ws.Range("Q7").Value = 0.148471615720524
$ws.Range("R7").Value = 0.07423580786026202
$ws.Range("S7").Value = 0.4279475982532751
$ws.Range("B8").Value = 0.1344339622641509
$ws.Range("D8").Value = 0.01650943396226415
$ws.Range("F8").Value = 0.06132075471698113
$ws.Range("J8").Value = 0.1202830188679245
$ws.Range("O8").Value = 0.02830188679245283
$ws.Range("Q8").Value = 0.1485849056603774
$ws.Range("R8").Value = 0.08726415094339622
$ws.Range("S8").Value = 0.4033018867924528
$ws.Range("B9").Value = 0.0949367088607595
$ws.Range("D9").Value = 0.03164556962025317
$ws.Range("F9").Value = 0.06329113924050633
$ws.Range("J9").Value = 0.1012658227848101
$ws.Range("O9").Value = 0.05063291139240506
$ws.Range("Q9").Value = 0.1835443037974684
$ws.Range("R9").Value = 0.06329113924050633
$ws.Range("S9").Value = 0.4113924050632912
$ws.Range("B10").Value = 0.1129900076863951
$ws.Range("D10").Value = 0.0207532667179093
$ws.Range("E10").Value = 0.0007686395080707148
$ws.Range("F10").Value = 0.07609531129900077
$ws.Range("J10").Value = 0.1137586471944658
$ws.Range("O10").Value = 0.02459646425826287
$ws.Range("Q10").Value = 0.2052267486548809
$ws.Range("R10").Value = 0.07993850883935434
$ws.Range("S10").Value = 0.3658724058416603
$ws.Range("F11").Value = 0.002976190476190476
$ws.Range("G11").Value = 0.1517857142857143
$ws.Range("J11").Value = 0.09523809523809523
$ws.Range("K11").Value = 0.2083333333333333
$ws.Range("L11").Value = 0.5327380952380952
$ws.Range("S11").Value = 0.008928571428571428
$ws.Range("G12").Value = 0.7663043478260869
$ws.Range("J12").Value = 0.1739130434782609
$ws.Range("L12").Value = 0.02717391304347826
$ws.Range("S12").Value = 0.03260869565217391
$ws.Range("G13").Value = 0.7894736842105263
$ws.Range("J13").Value = 0.2105263157894737
$ws.Range("F15").Value = 0.02380952380952381
$ws.Range("H15").Value = 0.1547619047619048
$ws.Range("I15").Value = 0.03571428571428571
$ws.Range("J15").Value = 0.3214285714285715
$ws.Range("K15").Value = 0.05952380952380952
$ws.Range("M15").Value = 0.003968253968253968
$ws.Range("N15").Value = 0.003968253968253968
$ws.Range("O15").Value = 0.05158730158730158
$ws.Range("S15").Value = 0.3452380952380952
$ws.Range("F16").Value = 0.01376146788990826
$ws.Range("H16").Value = 0.1467889908256881
$ws.Range("I16").Value = 0.07339449541284404
$ws.Range("J16").Value = 0.4495412844036697
$ws.Range("K16").Value = 0.1238532110091743
$ws.Range("M16").Value = 0.03669724770642202
$ws.Range("O16").Value = 0.03669724770642202
$ws.Range("S16").Value = 0.1192660550458716
$ws.Range("F17").Value = 0.02625298329355609
$ws.Range("H17").Value = 0.1670644391408115
$ws.Range("I17").Value = 0.09307875894988067
$ws.Range("J17").Value = 0.4081145584725537
$ws.Range("K17").Value = 0.1050119331742243
$ws.Range("M17").Value = 0.03579952267303103
$ws.Range("O17").Value = 0.05489260143198091
$ws.Range("S17").Value = 0.1097852028639618
$ws.Range("F18").Value = 0.02793296089385475
$ws.Range("H18").Value = 0.1229050279329609
$ws.Range("I18").Value = 0.07262569832402235
$ws.Range("J18").Value = 0.4636871508379888
$ws.Range("K18").Value = 0.1005586592178771
$ws.Range("M18").Value = 0.0223463687150838
$ws.Range("N18").Value = 0.00558659217877095
$ws.Range("O18").Value = 0.07262569832402235
$ws.Range("S18").Value = 0.111731843575419
$ws.Range("F19").Value = 0.02081727062451812
$ws.Range("H19").Value = 0.2012336160370085
$ws.Range("I19").Value = 0.06245181187355436
$ws.Range("J19").Value = 0.3747108712413261
$ws.Range("K19").Value = 0.1225905936777178
$ws.Range("M19").Value = 0.02313030069390902
$ws.Range("N19").Value = 0.0007710100231303007
$ws.Range("O19").Value = 0.07941403238242097
$ws.Range("S19").Value = 0.1148804934464148
